$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Percentage of estimates for test set that are off by less than 25% from true value: 76.00",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Percentage of estimates for test set that are off by less than 25% from true value: 73.54",
    2)

$d.Content.Find.Execute(
    "Percentage of estimates for test set that are off by less than 35% from true value: 84.31",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Percentage of estimates for test set that are off by less than 35% from true value: 85.23",
    2)

$d.Content.Find.Execute(
    "Percentage of estimates for test set that are off by less than 45% from true value: 93.54",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Percentage of estimates for test set that are off by less than 45% from true value: 92.62",
    2)

$d.Content.Find.Execute(
    "Percentage of estimates for test set that are off by less than 55% from true value: 96.00",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Percentage of estimates for test set that are off by less than 55% from true value: 96.92",
    2)
